# Actualiza la tabla de "Estado de Cuenta" (filas 16-21): se eliminan los
# registros previos y se reordenan/actualizan con los datos de la nueva
# base de datos (mismos trabajadores, nuevo orden y nuevo valor de mora).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C16").Value = "73199947"
$ws.Range("D16").Value = "GUSTAVO ADOLFO FELFLE FUENTES"
$ws.Range("E16").Value = "1711"
$ws.Range("G16").Value = 737717

$ws.Range("C17").Value = "9096389"
$ws.Range("D17").Value = "FERNANDO FRANCISCO FELFLE FUENTES"
$ws.Range("E17").Value = "1711"
$ws.Range("G17").Value = 737717

$ws.Range("C18").Value = "9096389"
$ws.Range("D18").Value = "FERNANDO FRANCISCO FELFLE FUENTES"
$ws.Range("E18").Value = "1710"
$ws.Range("G18").Value = 737717

$ws.Range("C19").Value = "9100677"
$ws.Range("D19").Value = "CARLOS ALBERTO HOYOS RIOS"
$ws.Range("E19").Value = "1711"
$ws.Range("G19").Value = 737717

$ws.Range("C20").Value = "9100677"
$ws.Range("D20").Value = "CARLOS ALBERTO HOYOS RIOS"
$ws.Range("E20").Value = "1710"
$ws.Range("G20").Value = 737717

$ws.Range("C21").Value = "1143393504"
$ws.Range("D21").Value = "OSCAR ENRIQUE ANAYA MEJIA"
$ws.Range("E21").Value = "1711"
$ws.Range("G21").Value = 737717
